$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A ("مواد مستقیم مصرفی" etc. row labels) gets a best-fit-style width;
# columns B:F (the yearly value columns) share a common best-fit-style width.
$ws.Columns.Item(1).ColumnWidth = 27.6
$ws.Range("B1:F1").EntireColumn.ColumnWidth = 19.6

# Restore the active cell/selection that was saved with the workbook (D14).
[void]$ws.Range("D14").Select()
